$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 0.01
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.99

# Row 5
$ws.Range("F5").Value = 0.01
$ws.Range("G5").Value = 0.99
$ws.Range("H5").Value = 0.01

# Row 6
$ws.Range("F6").Value = 0.01
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.99
